# Apply the "Corrected schematic for OPAMP connection" edit:
#  - Update the absolute-path hint recorded by Excel for this workbook
#  - Rename Sheet2 -> "PCB" and make it the active sheet
#  - Add a small "PCB correction log" table to the PCB sheet
#  - Move Sheet1's selection off the old header block

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- workbook-level bookkeeping -------------------------------------------------
# absPath mc:AlternateContent hint - the file now lives one level deeper
$wb.Path = "E:\Projects\Powersupply\PCB\PCB\Aimac\1.5KPS\"

# --- rename + repurpose Sheet2 as the PCB correction log ------------------------
$ws2.Name = "PCB"

# Header banner row (merged-looking single cell + two blank formatted neighbours)
$ws2.Range("B7").Value = "PCB VERSION -1 correction"
$ws2.Range("B7").Font.Bold = $true
$ws2.Range("B7").IndentLevel = 0
$ws2.Range("C7:D7").IndentLevel = 0

# Numbered correction list, rows 11-18
$corrections = @(
    "Resistor pad size reduce",
    "Resistor pad gap reduce",
    "cap pad size to 6mm",
    "New CT",
    "12V power pin and DIP RMC gap increse ",
    "NTC pad size increse",
    "NTC pad need to redesign( both pin not in stright line)"
)

for ($i = 0; $i -lt 8; $i++) {
    $row = 11 + $i
    $ws2.Cells.Item($row, 1).Value = $i + 1
    if ($i -lt $corrections.Length) {
        $ws2.Cells.Item($row, 2).Value = $corrections[$i]
    }
}
$ws2.Range("A11:A18").HorizontalAlignment = -4108
$ws2.Range("A11:A18").VerticalAlignment = -4108

# Column widths to match the new layout
$ws2.Columns.Item(2).ColumnWidth = 50.15

# --- sheet view / selection bookkeeping -----------------------------------------
$ws1.Range("C44").Select()

$ws2.Activate()
$ws2.Range("B21").Select()
